$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.140.90'
$ws.Range("E2").Value = '  -4.34%  '

# Row 3
$ws.Range("D3").Value = '1.650.95'
$ws.Range("E3").Value = '  -3.50%  '

# Row 5
$ws.Range("D5").Value = '''215.21'
$ws.Range("E5").Value = '  -4.09%  '

# Row 6
$ws.Range("D6").Value = '''0.5113'
$ws.Range("E6").Value = '  -3.19%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").Value = '''0.2585'
$ws.Range("E8").Value = '  -2.71%  '

# Row 9
$ws.Range("E9").Value = '  -3.71%  '

# Row 10
$ws.Range("E10").Value = '  -4.14%  '

# Row 11
$ws.Range("D11").Value = '''0.07788'
$ws.Range("E11").Value = '  +1.24%  '

# Row 12
$ws.Range("D12").Value = '1.655.27'
$ws.Range("E12").Value = '  -3.69%  '

# Row 13
$ws.Range("D13").Value = '''4.281'
$ws.Range("E13").Value = '  -4.86%  '

# Row 14
$ws.Range("D14").Value = '1.879.32'
$ws.Range("E14").Value = '  -3.44%  '

# Row 15
$ws.Range("D15").Value = '''0.5512'
$ws.Range("E15").Value = '  -5.75%  '

# Row 16
$ws.Range("D16").Value = '0.0₅8001'
$ws.Range("E16").Value = '  -2.72%  '

# Row 17
$ws.Range("D17").Value = '''64.05'
$ws.Range("E17").Value = '  -5.76%  '

# Row 18
$ws.Range("D18").Value = '26.141.95'
$ws.Range("E18").Value = '  -4.43%  '

# Row 19
$ws.Range("D19").Value = '''1.004'
$ws.Range("E19").Value = '  +0.05%  '

# Row 20
$ws.Range("D20").Value = '''210.96'
$ws.Range("E20").Value = '  -5.14%  '

# Row 21
$ws.Range("D21").Value = '''4.391'
$ws.Range("E21").Value = '  -5.29%  '

# Row 23
$ws.Range("D23").Value = '''6.040'
$ws.Range("E23").Value = '  +0.53%  '

# Row 24
$ws.Range("D24").Value = '''1.004'

# Row 25
$ws.Range("D25").Value = '''143.81'
$ws.Range("E25").Value = '  -0.59%  '

# Row 26
$ws.Range("D26").Value = '''1.750'
$ws.Range("E26").Value = '  +3.54%  '

# Row 27
$ws.Range("D27").Value = '''0.1174'
$ws.Range("E27").Value = '  -2.60%  '

# Row 28
$ws.Range("D28").Value = '''6.976'

# Row 29
$ws.Range("D29").Value = '''15.80'
$ws.Range("E29").Value = '  -2.53%  '

# Row 30
$ws.Range("D30").Value = '''0.05134'
$ws.Range("E30").Value = '  -3.65%  '

# Row 31
$ws.Range("D31").Value = '''1.241'
$ws.Range("E31").Value = '  -3.93%  '

# Row 32
$ws.Range("E32").Value = '  -3.49%  '

# Row 33
$ws.Range("D33").Value = '''3.215'
$ws.Range("E33").Value = '  -6.21%  '

# Row 34
$ws.Range("E34").Value = '  -4.76%  '

# Row 35
$ws.Range("D35").Value = '''2.738'
$ws.Range("E35").Value = '  -4.60%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''0.9236'
$ws.Range("E36").Value = '  -3.06%  '

# Row 37
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '''2.349'
$ws.Range("E37").Value = '  -1.85%  '

# Row 38
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.168.00'
$ws.Range("E38").Value = '  +1.53%  '

# Row 39
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '''0.5688'
$ws.Range("E39").Value = '  -2.54%  '

# Row 40
$ws.Range("D40").Value = '''0.01582'
$ws.Range("E40").Value = '  -3.23%  '

# Row 41
$ws.Range("D41").Value = '''2.551'
$ws.Range("E41").Value = '  -0.81%  '

# Row 42
$ws.Range("E42").Value = '  -0.02%  '

# Row 43
$ws.Range("E43").Value = '  -2.24%  '

# Row 44
$ws.Range("D44").Value = '''0.8240'
$ws.Range("E44").Value = '  -1.92%  '

# Row 45
$ws.Range("D45").Value = '''100.12'
$ws.Range("E45").Value = '  -1.14%  '

# Row 46
$ws.Range("D46").Value = '1.788.92'
$ws.Range("E46").Value = '  -3.48%  '

# Row 47
$ws.Range("D47").Value = '0.0₈117'
$ws.Range("E47").Value = '  +1.06%  '

# Row 48
$ws.Range("D48").Value = '''0.4548'
$ws.Range("E48").Value = '  -0.17%  '

# Row 49
$ws.Range("D49").Value = '''55.41'
$ws.Range("E49").Value = '  -3.80%  '

# Row 50
$ws.Range("E50").Value = '  +0.32%  '

# Row 51
$ws.Range("D51").Value = '''7.864'
$ws.Range("E51").Value = '  -2.99%  '
